$d = $word.ActiveDocument

# --- Change 1: add a new run " Лента заполнена пустыми символами." right
# after the run "Каретка зафиксирована и находится по центру." (its own
# paragraph), as a separate <w:r> (same run formatting as its neighbour).
$rng1 = $d.Content
$found = $rng1.Find.Execute("Каретка зафиксирована и находится по центру.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the caret paragraph text"
}
$rng1.Collapse(0)
$splitPos = $rng1.Start

$insertRange = $d.Range($splitPos, $splitPos)
$insertRange.InsertAfter(" Лента заполнена пустыми символами.")

# Force a run boundary at $splitPos (otherwise the newly inserted text
# gets merged into the preceding run, since both share identical rPr):
# add a throw-away bookmark exactly at the boundary, then remove it --
# the split itself survives the removal.
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("ZzTempSplit", $bmRange)
$d.Bookmarks("ZzTempSplit").Delete()

# --- Change 2: split the run " перемещается " into " пере" / "мещается "
# with a _GoBack bookmark in between (Word keeps only one _GoBack, so
# this also removes the stray one that used to sit at the end of the
# document).
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Лента перемещается", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'Лента перемещается'"
}
# $rng2 now covers exactly "Лента перемещается"; the split point is 10
# characters in ("Лента пере" = 10 chars), i.e. right before "мещается".
$splitPos2 = $rng2.Start + 10
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos2, $splitPos2))
